$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeout) values for rows 2-21, column G
$kValues = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 1
    6  = 0
    7  = 1
    8  = 2
    9  = 0
    10 = 1
    11 = 2
    12 = 0
    13 = 1
    14 = 1
    15 = 1
    16 = 0
    17 = 1
    18 = 2
    19 = 0
    20 = 2
    21 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
